$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the descriptive title by renaming the sheet tab to a generic name
$ws.Name = "Sheets1"

# Reset the view: scroll back to the top-left and move the selection to G8
$ws.Range("G8").Select()
